$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(2)
$tf = $shp.TextFrame
$tr = $tf.TextRange

$bullet = [char]0x2022

# --- Paragraph 5: "Position`t`t*  Source" -> "* Position`t`t * Source"
$para5 = $tr.Paragraphs(5)
$para5.Text = "$bullet Position`t`t $bullet Source"
$para5.ParagraphFormat.Bullet.Type = 0

# --- Paragraph 6: "BoundedBuffer`t*  CharUtil" (3 runs) gains a leading "* " run
#     and its middle run changes from "`t*  " to "`t *  " -> "`t * "
$para6 = $tr.Paragraphs(6)
$para6.ParagraphFormat.Bullet.Type = 0
$para6.InsertBefore("$bullet ")
# middle run is now run #2 (the tab/bullet run) - locate it via its text
$mid6 = $para6.Characters(1, $para6.Length)
# Replace "`t*  " with "`t *  " inside paragraph 6 text
$oldMid = "`t" + $bullet + "  "
$newMid = "`t " + $bullet + " "
$idx = $para6.Text.IndexOf($oldMid)
if ($idx -ge 0) {
    $run = $para6.Characters($idx + 1, $oldMid.Length)
    $run.Text = $newMid
}

# --- Paragraph 7: "ErrorHandler" gains a leading "* " run
$para7 = $tr.Paragraphs(7)
$para7.ParagraphFormat.Bullet.Type = 0
$para7.InsertBefore("$bullet ")

# --- Paragraph 9: "Symbol`t`t*  Token" -> "* Symbol`t`t * Token", gains Consolas font
$para9 = $tr.Paragraphs(9)
$para9.Text = "$bullet Symbol`t`t $bullet Token"
$para9.Font.Name = "Consolas"
$para9.ParagraphFormat.Bullet.Type = 0
